# Apply "last minute updates" to the first paragraph of the document:
#   1. Give the paragraph a box border (top/left/bottom/right, 5pt space to text).
#   2. Increase the paragraph's left indent from 120 (6pt) to 225 twips (11.25pt).
#   3. Rename the merge-field placeholder text and drop the now-redundant
#      trailing run that only held a single space.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# 1. Paragraph border on all four sides, 5pt (twips/20) space from the text.
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# 2. Left indent 120 -> 225 twips (Word COM uses points, so 225/20 = 11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# 3. Rename the placeholder id text in the first run only.
$rng = $p1.Range
$rng.Find.Execute("**ID__AFFARS_pgi_5317_topic_2__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5317__ID**", 2)

# Drop the trailing " " run left dangling at the end of the paragraph (before
# the paragraph mark) -- the author removed that whole run in the edit.
$pRange = $p1.Range
$text = $pRange.Text
$trimmed = $text.TrimEnd([char]13).TrimEnd(" ")
if ($trimmed.Length -lt $text.TrimEnd([char]13).Length) {
    $delStart = $pRange.Start + $trimmed.Length
    $delEnd = $pRange.End - 1
    $d.Range($delStart, $delEnd).Delete()
}
